# Weekly update for "Vega Modelo de Temuco - Espinaca" sheet.
# A new week of data is inserted at the top of the existing data block
# (rows 244-255), pushing the prior rows down by two rows
# (old 244-255 -> new 246-257), and the two freshly-inserted rows
# (244 and 245) are populated with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 244, shifting existing
# rows 244:255 down to 246:257 (formatting/styles carried along).
$ws.Rows.Item(244).Resize(2).Insert()

# ---- New row 244 ---------------------------------------------------
$ws.Range("A244").Value = 10
$ws.Range("B244").Value = "Vega Modelo de Temuco"
$ws.Range("C244").Value = "La Araucanía"
$ws.Range("D244").Value = 45008
$ws.Range("E244").Value = 9
$ws.Range("F244").Value = 100112012
$ws.Range("G244").Value = "Espinaca"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 35
$ws.Range("K244").Value = 10000
$ws.Range("L244").Value = 10000
$ws.Range("M244").Value = 10000
$ws.Range("N244").Value = "$/docena de atados"
$ws.Range("O244").Value = "Región de La Araucanía"
$ws.Range("P244").Value = 3333
$ws.Range("Q244").Value = 3
$ws.Range("R244").Value = "Hortaliza"

# ---- New row 245 -----------------------------------------------------
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 45008
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = 100112012
$ws.Range("G245").Value = "Espinaca"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Segunda"
$ws.Range("J245").Value = 30
$ws.Range("K245").Value = 9000
$ws.Range("L245").Value = 9000
$ws.Range("M245").Value = 9000
$ws.Range("N245").Value = "$/docena de atados"
$ws.Range("O245").Value = "Región de La Araucanía"
$ws.Range("P245").Value = 3000
$ws.Range("Q245").Value = 3
$ws.Range("R245").Value = "Hortaliza"
